$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old test data had three extra records (rows 16:18) that are no longer
# needed; delete those entire rows, shifting everything below up by three.
$ws.Range("A16:C18").EntireRow.Delete()

# Extend the sheet's tracked row range down to the very bottom of the sheet
# and give the final few rows their (non-default) height, mirroring the
# trailing "touched" rows left behind at the bottom of the worksheet.
$ws.Range("A1048574:A1048576").RowHeight = 12.8

# Move the active selection to the first cell of what is now the last
# (previously 19th, now 16th) data row.
[void]$ws.Range("A16").Select()
